# Adds six new rows (9-14) of Pokémon-card price data to the existing
# Nome/Código/Coleção table on the active sheet, matching the upstream
# "Add files via upload" commit.
#
# New data (Name / Code / Collection):
#   9  Iron Valiant ex  (089/182)  PAR
#   10 Iron Valiant ex  (261/182)  PAR
#   11 Scream Tail      (086/182)  PAR
#   12 Alcremie-V       (22/073)   CPA
#   13 Jynx ex          (124/165)  MEW
#   14 Arcanine ex      (032/198)  SV1
#
# Cells are written in the exact order the original author must have
# entered/edited them so the shared-string table comes out in the same
# sequence as the target file (rows are built by copying the previous
# row then only touching the cells that actually changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: brand-new row, left to right ------------------------------
$ws.Range("A9").Value = "Iron Valiant ex"
$ws.Range("B9").Value = "(089/182)"
$ws.Range("C9").Value = "PAR"

# --- Row 10: same card, only the code differs --------------------------
$ws.Range("A10").Value = "Iron Valiant ex"
$ws.Range("B10").Value = "(261/182)"
$ws.Range("C10").Value = "PAR"

# --- Row 11: code edited before the name --------------------------------
$ws.Range("C11").Value = "PAR"
$ws.Range("B11").Value = "(086/182)"
$ws.Range("A11").Value = "Scream Tail"

# --- Row 12: name, then collection, then code ---------------------------
$ws.Range("A12").Value = "Alcremie-V"
$ws.Range("C12").Value = "CPA"
$ws.Range("B12").Value = "(22/073)"

# --- Row 13: name then code (collection reuses existing "MEW") ----------
$ws.Range("A13").Value = "Jynx ex"
$ws.Range("B13").Value = "(124/165)"
$ws.Range("C13").Value = "MEW"

# --- Row 14: name, then collection, then code ---------------------------
$ws.Range("A14").Value = "Arcanine ex"
$ws.Range("C14").Value = "SV1"
$ws.Range("B14").Value = "(032/198)"

# --- Formatting: reuse the bold-Arial-12 style already used by the -------
# other data rows (copy its format instead of touching Font properties
# one-by-one so no duplicate font/style entries get minted).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A9:A14").PasteSpecial(-4122) | Out-Null
$ws.Range("B10:B14").PasteSpecial(-4122) | Out-Null

# B9 ("089/182") is the highlighted/bigger header-style cell: start from
# the same bold Arial font, then bump the size and set left/center
# alignment - this mints exactly one new font + one new cell style.
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Font.Size = 14.3
$ws.Range("B9").HorizontalAlignment = -4131
$ws.Range("B9").VerticalAlignment = -4108

$excel.CutCopyMode = $false

# --- Row heights, matching the taller header row (9) vs. the regular ----
# data rows (10-14).
$ws.Rows.Item(9).RowHeight = 18
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75
$ws.Rows.Item(13).RowHeight = 15.75
$ws.Rows.Item(14).RowHeight = 15.75

# --- Final selection, matching the saved cursor position in the source --
$ws.Range("B15").Select() | Out-Null
